$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Merge "Diba" + " Shojaeigoradel" runs into a single run and drop the
#        stray _GoBack bookmark that used to sit between them. ---
$nameCell = $t.Cell(1, 3)
$nameCell.Range.Find.Execute("Diba Shojaeigoradel", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Diba Shojaeigoradel", 2) | Out-Null

# --- 2. Fill in the score "5" for the middle ("Umar Ehsan") column on the
#        first four data rows. ---
for ($row = 2; $row -le 5; $row++) {
    $scoreCell = $t.Cell($row, 2)
    $scoreCell.Range.Text = "5"
}

# --- 3. Fifth (last) data row: type "5" plus a throwaway trailing marker so
#        the bookmark insertion point isn't sitting exactly on the
#        paragraph-end boundary, add the _GoBack bookmark right after the
#        "5", then strip the marker back out again. ---
$lastCell = $t.Cell(6, 2)
$lastCell.Range.Text = "5X"
$afterFive = $lastCell.Range.Start + 1
$pt = $d.Range($afterFive, $afterFive)
$d.Bookmarks.Add("_GoBack", $pt)

$marker = $d.Range($afterFive, $afterFive + 1)
$marker.Text = ""
